# Redefine efficiency of electrolysis
#
# The "tech_data" sheet stores H2 electrolysis efficiency (col I, ACT_EFF) as a
# fuel-input-per-unit-of-output figure (0.665). The model now wants the
# reciprocal (output per unit input, ~1.504) so the base-year cells become a
# formula "=1/0.665" and the escalation formulas that reference them
# (I*1.04, I*1.14) recompute accordingly. The new number format ("0.00") is
# also applied to the previously-blank interpolation years in column I so the
# whole block is formatted consistently, and the stray top-border formatting
# that had been left on F17/F18 is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tech_data")

# --- Electroliser small (1 MW) block (rows 4-9) ---
$ws.Range("I4").Formula = "=1/0.665"
$ws.Range("I4").NumberFormat = "0.00"

$ws.Range("I5").NumberFormat = "0.00"

$ws.Range("I6").NumberFormat = "0.00"          # formula =I4*1.04 stays, recalculates
$ws.Range("I7").NumberFormat = "0.00"
$ws.Range("I8").NumberFormat = "0.00"

$ws.Range("I9").NumberFormat = "0.00"          # formula =I6*1.14 stays, recalculates

# --- Electroliser large (100 MW) block (rows 10-15) ---
$ws.Range("I10").Formula = "=1/0.665"
$ws.Range("I10").NumberFormat = "0.00"

$ws.Range("I11").NumberFormat = "0.00"

$ws.Range("I12").NumberFormat = "0.00"         # formula =I10*1.04 stays, recalculates
$ws.Range("I13").NumberFormat = "0.00"
$ws.Range("I14").NumberFormat = "0.00"

$ws.Range("I15").NumberFormat = "0.00"         # formula =I12*1.14 stays, recalculates

# --- Drop the leftover top border on the compression rows ---
$ws.Range("F17").Borders.Item(8).LineStyle = -4142
$ws.Range("F18").Borders.Item(8).LineStyle = -4142

# --- Restore the cursor to where the author left it ---
$ws.Range("K22").Select()
